$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='56.360.92'}
    @{Cell='E2'; Value='  +9.41%  '}
    @{Cell='D3'; Value='3.234.35'}
    @{Cell='E3'; Value='  +4.25%  '}
    @{Cell='E4'; Value='  -0.09%  '}
    @{Cell='D5'; Value='398.93'}
    @{Cell='E5'; Value='  +3.35%  '}
    @{Cell='D6'; Value='110.99'}
    @{Cell='E6'; Value='  +6.59%  '}
    @{Cell='E8'; Value='  -0.05%  '}
    @{Cell='E9'; Value='  +6.36%  '}
    @{Cell='D10'; Value='39.51'}
    @{Cell='E10'; Value='  +6.04%  '}
    @{Cell='D11'; Value='0.0907'}
    @{Cell='E11'; Value='  +5.69%  '}
    @{Cell='E12'; Value='  +2.11%  '}
    @{Cell='D13'; Value='3.739.78'}
    @{Cell='E13'; Value='  +4.05%  '}
    @{Cell='D14'; Value='8.12'}
    @{Cell='E14'; Value='  +3.64%  '}
    @{Cell='D15'; Value='19.09'}
    @{Cell='E15'; Value='  +2.75%  '}
    @{Cell='D16'; Value='3.217.71'}
    @{Cell='E16'; Value='  +3.75%  '}
    @{Cell='E17'; Value='  +5.57%  '}
    @{Cell='D18'; Value='10.68'}
    @{Cell='E18'; Value='  -2.40%  '}
    @{Cell='D19'; Value='56.136.96'}
    @{Cell='D20'; Value='3.34'}
    @{Cell='E20'; Value='  +1.87%  '}
    @{Cell='E21'; Value='  +5.82%  '}
    @{Cell='D22'; Value='13.07'}
    @{Cell='E22'; Value='  +4.14%  '}
    @{Cell='D23'; Value='304.41'}
    @{Cell='D24'; Value='75.49'}
    @{Cell='E24'; Value='  +7.73%  '}
    @{Cell='D25'; Value='3.23'}
    @{Cell='E25'; Value='  +1.74%  '}
    @{Cell='D26'; Value='8.20'}
    @{Cell='E26'; Value='  +1.10%  '}
    @{Cell='D27'; Value='28.29'}
    @{Cell='E27'; Value='  +2.99%  '}
    @{Cell='D28'; Value='7.47'}
    @{Cell='E28'; Value='  +3.36%  '}
    @{Cell='D29'; Value='0.174'}
    @{Cell='E29'; Value='  +4.18%  '}
    @{Cell='E30'; Value='  -0.10%  '}
    @{Cell='E31'; Value='  +4.23%  '}
    @{Cell='D32'; Value='11.19'}
    @{Cell='E32'; Value='  +7.31%  '}
    @{Cell='D33'; Value='0.0493'}
    @{Cell='E33'; Value='  +3.14%  '}
    @{Cell='D34'; Value='36.42'}
    @{Cell='E34'; Value='  +2.06%  '}
    @{Cell='E35'; Value='  +3.51%  '}
    @{Cell='D36'; Value='51.40'}
    @{Cell='E36'; Value='  +2.51%  '}
    @{Cell='D37'; Value='3.14'}
    @{Cell='E37'; Value='  +24.71%  '}
    @{Cell='B38'; Value='LidoDAOToken'}
    @{Cell='C38'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'}
    @{Cell='D38'; Value='3.52'}
    @{Cell='E38'; Value='  +4.05%  '}
    @{Cell='B39'; Value='FirstDigitalUSD'}
    @{Cell='C39'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'}
    @{Cell='D39'; Value='1.00'}
    @{Cell='E39'; Value='  +0.07%  '}
    @{Cell='D40'; Value='134.82'}
    @{Cell='E40'; Value='  +4.49%  '}
    @{Cell='E41'; Value='  +3.10%  '}
    @{Cell='E42'; Value='  +7.01%  '}
    @{Cell='D43'; Value='17.23'}
    @{Cell='E43'; Value='  +3.31%  '}
    @{Cell='E44'; Value='  +3.24%  '}
    @{Cell='D45'; Value='0.284'}
    @{Cell='E45'; Value='  -2.89%  '}
    @{Cell='D46'; Value='22.39'}
    @{Cell='E46'; Value='  +0.59%  '}
    @{Cell='B47'; Value='ThetaToken'}
    @{Cell='C47'; Value='https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'}
    @{Cell='E47'; Value='  +46.76%  '}
    @{Cell='B48'; Value='WEMIXToken'}
    @{Cell='C48'; Value='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'}
    @{Cell='D48'; Value='2.12'}
    @{Cell='E48'; Value='  +1.97%  '}
    @{Cell='D49'; Value='2.49'}
    @{Cell='E49'; Value='  -1.33%  '}
    @{Cell='D50'; Value='2.140.78'}
    @{Cell='E50'; Value='  +3.06%  '}
    @{Cell='D51'; Value='0.0363'}
    @{Cell='E51'; Value='  +8.41%  '}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.Value = "'" + $u.Value
    $r.ClearFormats()
}
